$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Helpers
# ----------------------------------------------------------------------

function Split-RunAt($pos) {
    # Force a run boundary at absolute character position $pos by toggling
    # bold on/off for the single character that starts right after $pos.
    # The net formatting is unchanged (Bold ends up unset, exactly as
    # before) but the engine must re-serialize, splitting the run from its
    # neighbour without disturbing the inherited <w:rPr> (e.g. the
    # <w:lang w:val="en-US"/> every run in this document carries).
    $r = $d.Range($pos, $pos + 1)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

function Find-Phrase($phrase) {
    $rng = $d.Content
    $f = $rng.Find
    $f.ClearFormatting()
    $f.Text = $phrase
    $f.Execute() | Out-Null
    if (-not $f.Found) { throw "phrase not found: $phrase" }
    return $rng
}

function Insert-EmptyParagraphAfterPhrase($phrase) {
    # Appends a brand-new, completely empty BodyText paragraph right after
    # the paragraph containing $phrase, using the "^p" Find/Replace
    # wildcard. This is the one technique that yields a clean <w:p> with
    # no stray run inside (matching genuinely-empty paragraphs already in
    # this document), unlike Range.InsertParagraphAfter()/Before() which
    # always synthesizes a spurious empty <w:r>.
    $rng = $d.Content
    $f = $rng.Find
    $f.ClearFormatting()
    $f.Execute($phrase, $true, $false, $false, $false, $false, $true, 1, $false,
               "$phrase^p", 2) | Out-Null
    if (-not $f.Found) { throw "phrase not found for empty-paragraph insert: $phrase" }
}

function Append-ParagraphWithText($afterPhrase, $newText) {
    # Appends $newText as a brand-new paragraph right after the paragraph
    # containing $afterPhrase. To make sure the new paragraph's runs carry
    # the same <w:rPr> (<w:lang w:val="en-US"/>) as the rest of the
    # document, the text is first appended *inline* (same paragraph, so it
    # inherits the formatting of $afterPhrase's run), then promoted to its
    # own paragraph by inserting a paragraph break exactly at the old/new
    # text boundary - that boundary sits inside a run that already carries
    # the correct <w:rPr>, so the split keeps it (same trick used for
    # Split-RunAt).
    $rng = $d.Content
    $f = $rng.Find
    $f.ClearFormatting()
    $f.Execute($afterPhrase, $true, $false, $false, $false, $false, $true, 1, $false,
               "$afterPhrase$newText", 2) | Out-Null
    if (-not $f.Found) { throw "phrase not found for paragraph append: $afterPhrase" }

    $rng2 = $d.Content
    $f2 = $rng2.Find
    $f2.ClearFormatting()
    $f2.Text = $newText
    $f2.Execute() | Out-Null
    if (-not $f2.Found) { throw "appended text not found: $newText" }
    $splitPos = $rng2.Start
    $zeroRng = $d.Range($splitPos, $splitPos)
    $zeroRng.InsertParagraphBefore()
}

# ---------------------------------------------------------------------
# Change 1: "...get it to work and I continued on to the first lesson..."
#        -> "...get it to work, and I continued the first lesson..."
#           split across 5 runs, plus a new empty paragraph right after.
# ---------------------------------------------------------------------
$rng = Find-Phrase "work and I continued on to the first lesson"
$matchStart = $rng.Start
$rng.Text = "work, and I continued the first lesson"

$seg1 = "work,"
$seg2 = " and I "
$seg3 = "continued"

$pos = $matchStart + $seg1.Length
Split-RunAt $pos
$pos = $pos + $seg2.Length
Split-RunAt $pos
$pos = $pos + $seg3.Length
Split-RunAt $pos

Insert-EmptyParagraphAfterPhrase "the first lesson: a crash course in node."

# ---------------------------------------------------------------------
# Change 2: "...install as well which I found out..."
#        -> "...install as well, which I found out..." split into 3 runs.
# ---------------------------------------------------------------------
$rng = Find-Phrase "as well which I found out"
$matchStart = $rng.Start
$rng.Text = "as well, which I found out"

$seg1 = "as well, which"
$pos = $matchStart + $seg1.Length
Split-RunAt $pos

# ---------------------------------------------------------------------
# Change 3: insert a new empty BodyText paragraph before the paragraph
# that begins "Following the MongoDB lesson..." i.e. right after the
# "...Stack exchange." paragraph.
# ---------------------------------------------------------------------
Insert-EmptyParagraphAfterPhrase "I found out after a brief visit to Stack exchange."

# ---------------------------------------------------------------------
# Change 4: "...MongoDB cheat sheet I downloaded..." ('sheet' wrapped in
# proofErr gramStart/gramEnd) -> "...MongoDB cheat sheet, I downloaded..."
# (comma added; the grammar-error marker is gone because 'sheet' is no
# longer directly followed by 'I' with no punctuation).
# ---------------------------------------------------------------------
$rng = Find-Phrase "cheat sheet I downloaded"
$matchStart = $rng.Start
$rng.Text = "cheat sheet, I downloaded"

$seg1 = "cheat sheet,"
$pos = $matchStart + $seg1.Length
Split-RunAt $pos

# ---------------------------------------------------------------------
# Change 5: append new diary content after "...for later reference use."
#   - empty paragraph
#   - "Next I checked out the Express lesson..." paragraph (several runs)
#   - empty paragraph
#   - "The new lesson did not have any proper error handling..." paragraph
# ---------------------------------------------------------------------
Insert-EmptyParagraphAfterPhrase "for later reference use."

$expressText = "Next I checked out the Express lesson, and even though I have setup an express server before, I learned quite a lot of new things like the server rendering with handlebars and the proper use of middleware, which always was slightly confusing before. Postman was already familiar to me and using it to test the API was pretty straightforward."
Append-ParagraphWithText "for later reference use." $expressText

# split the Express paragraph into the multiple runs shown by the diff
$rng = Find-Phrase $expressText
$matchStart = $rng.Start
$parts = @(
    "Next I checked out the Express lesson, and even though I have setup an express server before, I learned quite a lot of new things like the server rendering with handlebars and the proper ",
    "use of middleware, which always was ",
    "slightly ",
    "confusing ",
    "before. ",
    "Postman was already familiar to me and using it to test the ",
    "API was pretty straightforward."
)
$pos = $matchStart
for ($i = 0; $i -lt $parts.Length - 1; $i++) {
    $pos = $pos + $parts[$i].Length
    Split-RunAt $pos
}

Insert-EmptyParagraphAfterPhrase "API was pretty straightforward."

$finalText = "The new lesson did not have any proper error handling, but I installed the express validator package without many issues and added the validator code to the tutorial express setup. While the video did not touch any further on the passport package it is something I will check out in more detail later and possibly use in my own project."
Append-ParagraphWithText "API was pretty straightforward." $finalText

$rng = Find-Phrase $finalText
$matchStart = $rng.Start
$parts2 = @(
    "The new lesson did not have any proper error ",
    "handling,",
    " but I installed the ",
    "express validator package without ",
    "many",
    " issues and added the validator code to the ",
    "tutorial express setup. While the video did not touch any further on the passport package it is something I will check out in more detail later and possibly use in my own project."
)
$pos = $matchStart
for ($i = 0; $i -lt $parts2.Length - 1; $i++) {
    $pos = $pos + $parts2[$i].Length
    Split-RunAt $pos
}

Write-Host "done"
